# Apply the latest cryptos-list scrape: refreshed prices/volumes for most
# rows, plus Binance-PegBSC-USD (row 29) and Bittensor (row 30) swapping
# rank position (now Bittensor #29 / Binance-PegBSC-USD #30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='67.175.06'},
    @{Cell='E2'; Value='  -0.70%  '},
    @{Cell='D3'; Value='2.606.37'},
    @{Cell='E3'; Value='  -0.43%  '},
    @{Cell='E4'; Value='  +0.01%  '},
    @{Cell='D5'; Value='''590.86'},
    @{Cell='E5'; Value='  -1.89%  '},
    @{Cell='D6'; Value='''150.44'},
    @{Cell='E6'; Value='  -2.49%  '},
    @{Cell='D8'; Value='''0.546'},
    @{Cell='E8'; Value='  -0.67%  '},
    @{Cell='D9'; Value='2.604.55'},
    @{Cell='E9'; Value='  -0.37%  '},
    @{Cell='E10'; Value='  +0.73%  '},
    @{Cell='E11'; Value='  -0.09%  '},
    @{Cell='D12'; Value='''5.16'},
    @{Cell='E12'; Value='  -1.57%  '},
    @{Cell='D13'; Value='''0.343'},
    @{Cell='E13'; Value='  -3.01%  '},
    @{Cell='D14'; Value='''27.23'},
    @{Cell='E14'; Value='  -2.62%  '},
    @{Cell='D15'; Value='3.080.39'},
    @{Cell='E15'; Value='  -0.43%  '},
    @{Cell='E16'; Value='  -3.30%  '},
    @{Cell='D17'; Value='67.033.84'},
    @{Cell='E17'; Value='  -0.84%  '},
    @{Cell='D18'; Value='2.601.70'},
    @{Cell='E18'; Value='  -0.49%  '},
    @{Cell='D19'; Value='''374.79'},
    @{Cell='E19'; Value='  +3.18%  '},
    @{Cell='D20'; Value='''11.04'},
    @{Cell='E20'; Value='  -2.00%  '},
    @{Cell='D21'; Value='''7.38'},
    @{Cell='E21'; Value='  -3.16%  '},
    @{Cell='E22'; Value='  -0.60%  '},
    @{Cell='D23'; Value='''4.75'},
    @{Cell='E23'; Value='  -4.64%  '},
    @{Cell='E24'; Value='  -4.68%  '},
    @{Cell='D25'; Value='''73.26'},
    @{Cell='E25'; Value='  +4.59%  '},
    @{Cell='E26'; Value='  -0.06%  '},
    @{Cell='D27'; Value='''9.94'},
    @{Cell='E27'; Value='  -1.20%  '},
    @{Cell='D28'; Value='2.738.12'},
    @{Cell='E28'; Value='  -0.24%  '},
    @{Cell='B29'; Value='Bittensor'},
    @{Cell='C29'; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'},
    @{Cell='D29'; Value='''582.94'},
    @{Cell='E29'; Value='  -0.56%  '},
    @{Cell='B30'; Value='Binance-PegBSC-USD'},
    @{Cell='C30'; Value='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'},
    @{Cell='D30'; Value='''1.00'},
    @{Cell='E30'; Value='  -3.34%  '},
    @{Cell='D31'; Value='0.0₃0987'},
    @{Cell='E31'; Value='  -6.53%  '},
    @{Cell='E32'; Value='  -5.55%  '},
    @{Cell='D33'; Value='''7.67'},
    @{Cell='E33'; Value='  -3.50%  '},
    @{Cell='D34'; Value='''1.81'},
    @{Cell='E34'; Value='  -3.37%  '},
    @{Cell='E35'; Value='  +0.01%  '},
    @{Cell='E36'; Value='  -3.85%  '},
    @{Cell='E37'; Value='  -3.20%  '},
    @{Cell='D38'; Value='''156.57'},
    @{Cell='E38'; Value='  -0.03%  '},
    @{Cell='D39'; Value='''19.05'},
    @{Cell='E39'; Value='  -2.01%  '},
    @{Cell='D40'; Value='''0.365'},
    @{Cell='E40'; Value='  -1.78%  '},
    @{Cell='D41'; Value='''1.85'},
    @{Cell='E41'; Value='  -0.61%  '},
    @{Cell='D42'; Value='''5.24'},
    @{Cell='E42'; Value='  -3.40%  '},
    @{Cell='E43'; Value='  -4.18%  '},
    @{Cell='D44'; Value='''17.11'},
    @{Cell='E44'; Value='  +4.20%  '},
    @{Cell='E45'; Value='  -0.05%  '},
    @{Cell='D46'; Value='''153.22'},
    @{Cell='E46'; Value='  -2.44%  '},
    @{Cell='D47'; Value='0.0₆0285'},
    @{Cell='E47'; Value='  -1.91%  '},
    @{Cell='E48'; Value='  -1.47%  '},
    @{Cell='E49'; Value='  -1.69%  '},
    @{Cell='E50'; Value='  -4.43%  '},
    @{Cell='D51'; Value='''21.32'},
    @{Cell='E51'; Value='  +1.24%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
